# Insert a new data row at row 253 (pushing the existing rows 253:314 down
# to 254:315) and populate it with the new record's values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(253).Insert()

$ws.Cells.Item(253, 1).Value  = 7
$ws.Cells.Item(253, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(253, 3).Value  = "Ñuble"
$ws.Cells.Item(253, 4).Value  = 44508
$ws.Cells.Item(253, 5).Value  = 16
$ws.Cells.Item(253, 6).Value  = 100112020
$ws.Cells.Item(253, 7).Value  = "Tomate"
$ws.Cells.Item(253, 8).Value  = "Larga vida"
$ws.Cells.Item(253, 9).Value  = "Primera"
$ws.Cells.Item(253, 10).Value = 240
$ws.Cells.Item(253, 11).Value = 7500
$ws.Cells.Item(253, 12).Value = 8000
$ws.Cells.Item(253, 13).Value = 7750
$ws.Cells.Item(253, 14).Value = "`$/caja 10 kilos"
$ws.Cells.Item(253, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(253, 16).Value = 775
$ws.Cells.Item(253, 17).Value = 10
$ws.Cells.Item(253, 18).Value = "Hortaliza"
